# virtualan_collection_testcase_6.xlsx — "Updated with Field Name with standard
# name and Minor release"
#
# Renames several header cells in row 1 of the "API-Testing" sheet to their
# new standardized field names, normalizes the numeric value in L2, and
# updates the active selection, matching the commit's field-name cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) field-name standardization -------------------------
$ws.Range("E1").Value = "FormParams"        # was RequestParams
$ws.Range("G1").Value = "RequestHeaders"    # was RequestProcessingType
$ws.Range("J1").Value = "Action"            # was HTTPAction
$ws.Range("K1").Value = "ExcludeFields"     # was ExcludeField
$ws.Range("L1").Value = "StatusCode"        # was HttpStatusCode
$ws.Range("O1").Value = "Security"          # was security
$ws.Range("P1").Value = "Tags"              # was tags

# --- Minor-release numeric cleanup -----------------------------------------
$ws.Range("L2").Value = 415

# --- Selection / active cell ------------------------------------------------
$ws.Range("Q1").Select()
